$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "BOL"

$ws2.Range("A1").Value = "OrderId"
$ws2.Range("B1").Value = "Shipment Type"
$ws2.Range("C1").Value = "Starting tracking number"
$ws2.Range("D1").Value = "Order Status"

$ws2.Range("A2").Value = "51488767"
$ws2.Range("B2").Value = "Parcel"
$ws2.Range("C2").Value = "1z"

$ws2.Range("A3").Value = "51488755"
$ws2.Range("B3").Value = "LTL Order"
$ws2.Range("C3").Value = "fces"
$ws2.Range("D3").Value = "Confirm"

$ws2.Range("A4").Value = "51488755"
$ws2.Range("B4").Value = "LTL Order"
$ws2.Range("C4").Value = "fces"
$ws2.Range("D4").Value = "Withdraw"
